$wb = $excel.ActiveWorkbook

# Rename the "GET_Tests" sheet to "GET Tests"
$getSheet = $wb.Worksheets.Item("GET_Tests")
$getSheet.Name = "GET Tests"

# Add a new "Auth Tests" worksheet after the last existing sheet (POST Tests),
# so it lands at the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$authSheet = $wb.Worksheets.Add($null, $lastSheet)
$authSheet.Name = "Auth Tests"

# Make the new sheet the active/selected tab
$authSheet.Activate()
